# Updates the cryptos worksheet cell values per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.343.08"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.58"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.38"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  +7.09%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.62"
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.346"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.16"
$ws.Range("E10").Value = "  +13.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0720"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.185.18"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.35"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.699"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.912.90"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.366.69"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.20"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0822"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.92"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +27.23%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.02"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.34"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.151.16"
$ws.Range("E31").Value = "  +21.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.15"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0565"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +13.47%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.10"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0650"
$ws.Range("E41").Value = "  +8.55%  "
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.25"
$ws.Range("E43").Value = "  +6.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "90.09"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.339.51"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("B46").Value = "MultiversX"
$ws.Range("C46").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.93"
$ws.Range("E46").Value = "  +39.40%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.096.10"
$ws.Range("E51").Value = "  +1.80%  "
